# Weekly refresh of "Hortaliza, Terminal Hortofrutícola Agro Chillán - Lechuga"
# Two new price report rows (for 2022-07-27) are inserted at the top of the
# data set (rows 655-656), pushing the existing rows down by two and
# extending the table from 715 to 717 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 655:656 - everything currently at row 655
# downward (through 715) shifts down to 657..717. The inserted rows inherit
# the surrounding rows' formatting (e.g. the date number format on column D).
$ws.Rows("655:656").Insert()

# New row 655: Lechuga "Conconina(o)"
$ws.Cells.Item(655, 1).Value = 7
$ws.Cells.Item(655, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(655, 3).Value = "Ñuble"
$ws.Cells.Item(655, 4).Value = 44769
$ws.Cells.Item(655, 5).Value = 16
$ws.Cells.Item(655, 6).Value = 100112033
$ws.Cells.Item(655, 7).Value = "Lechuga"
$ws.Cells.Item(655, 8).Value = "Conconina(o)"
$ws.Cells.Item(655, 9).Value = "Primera"
$ws.Cells.Item(655, 10).Value = 160
$ws.Cells.Item(655, 11).Value = 7000
$ws.Cells.Item(655, 12).Value = 7500
$ws.Cells.Item(655, 13).Value = 7250
$ws.Cells.Item(655, 14).Value = "`$/caja 10 unidades"
$ws.Cells.Item(655, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(655, 16).Value = 725
$ws.Cells.Item(655, 17).Value = 10
$ws.Cells.Item(655, 18).Value = "Hortaliza"

# New row 656: Lechuga "Escarola"
$ws.Cells.Item(656, 1).Value = 7
$ws.Cells.Item(656, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(656, 3).Value = "Ñuble"
$ws.Cells.Item(656, 4).Value = 44769
$ws.Cells.Item(656, 5).Value = 16
$ws.Cells.Item(656, 6).Value = 100112033
$ws.Cells.Item(656, 7).Value = "Lechuga"
$ws.Cells.Item(656, 8).Value = "Escarola"
$ws.Cells.Item(656, 9).Value = "Primera"
$ws.Cells.Item(656, 10).Value = 120
$ws.Cells.Item(656, 11).Value = 9500
$ws.Cells.Item(656, 12).Value = 10000
$ws.Cells.Item(656, 13).Value = 9750
$ws.Cells.Item(656, 14).Value = "`$/caja 15 unidades"
$ws.Cells.Item(656, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(656, 16).Value = 650
$ws.Cells.Item(656, 17).Value = 15
$ws.Cells.Item(656, 18).Value = "Hortaliza"
